# Auto-generated edit script: applies the diff described in the task.
# Pattern: rows 3 & 4 swap/update record content (taxon "Garnlav" <-> "Granticka"),
# rows 25/26/27 rotate record content, and many other rows just bump the
# "Taxonsorteringsordning" (B column) sort-order value by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 79244

# Row 3
$ws.Range("A3").Value = 130862980
$ws.Range("B3").Value = 79244
$ws.Range("E3").Value = 6425
$ws.Range("F3").Value = 'Garnlav'
$ws.Range("G3").Value = 'Alectoria sarmentosa'
$ws.Range("H3").Value = '(Ach.) Ach.'
$ws.Range("I3").Value = 'Söderåsen Väst, Jmt'
$ws.Range("N3").Value = ""
$ws.Range("Q3").Value = 447218
$ws.Range("R3").Value = 7042948
$ws.Range("AC3").Value = 'På en tydligt gammal gran med nedåthängande grenar.'
$ws.Range("AF3").Value = ""
$ws.Range("AG3").Value = $False
$ws.Range("AH3").Value = ""
$ws.Range("AJ3").Value = ""
$ws.Range("AK3").Value = ""
$ws.Range("AM3").Value = ""
$ws.Range("AO3").Value = ""

# Row 4
$ws.Range("A4").Value = 130862972
$ws.Range("B4").Value = 91829
$ws.Range("E4").Value = 5432
$ws.Range("F4").Value = 'Granticka'
$ws.Range("G4").Value = 'Porodaedalea chrysoloma s.lat.'
$ws.Range("H4").Value = 'teleomorf'
$ws.Range("I4").Value = ""
$ws.Range("N4").Value = 'Söderåsen Väst, Jmt'
$ws.Range("Q4").Value = 447167
$ws.Range("R4").Value = 7042999
$ws.Range("AC4").Value = 'Flera fruktkroppar i en gammal relativt grov döende stående gran med full längd och 36 cm i brösthöjdsdiameter..'
$ws.Range("AF4").Value = "'0"
$ws.Range("AG4").Value = ""
$ws.Range("AH4").Value = 'Granskog'
$ws.Range("AJ4").Value = 'gran'
$ws.Range("AK4").Value = 'Picea abies'
$ws.Range("AM4").Value = 'Trädstam på levande träd'
$ws.Range("AO4").Value = 'Stem on living tree # Picea abies'

# Row 5
$ws.Range("B5").Value = 79244

# Row 6
$ws.Range("B6").Value = 79244

# Row 7
$ws.Range("B7").Value = 79244

# Row 8
$ws.Range("B8").Value = 79244

# Row 9
$ws.Range("B9").Value = 79244

# Row 10
$ws.Range("B10").Value = 79244

# Row 11
$ws.Range("B11").Value = 79244

# Row 12
$ws.Range("B12").Value = 79244

# Row 13
$ws.Range("B13").Value = 79244

# Row 14
$ws.Range("B14").Value = 79244

# Row 15
$ws.Range("B15").Value = 91829

# Row 16
$ws.Range("B16").Value = 79244

# Row 17
$ws.Range("B17").Value = 79244

# Row 18
$ws.Range("B18").Value = 79244

# Row 19
$ws.Range("B19").Value = 79244

# Row 20
$ws.Range("B20").Value = 79244

# Row 21
$ws.Range("B21").Value = 79244

# Row 22
$ws.Range("B22").Value = 91809

# Row 23
$ws.Range("B23").Value = 79244

# Row 24
$ws.Range("B24").Value = 79244

# Row 25
$ws.Range("A25").Value = 130865712
$ws.Range("B25").Value = 79244
$ws.Range("E25").Value = 6425
$ws.Range("F25").Value = 'Garnlav'
$ws.Range("G25").Value = 'Alectoria sarmentosa'
$ws.Range("H25").Value = '(Ach.) Ach.'
$ws.Range("Q25").Value = 447165
$ws.Range("R25").Value = 7043032
$ws.Range("S25").Value = 10
$ws.Range("Z25").Value = '15:13'
$ws.Range("AB25").Value = '15:13'

# Row 26
$ws.Range("A26").Value = 130865713
$ws.Range("B26").Value = 91805
$ws.Range("E26").Value = 1108
$ws.Range("F26").Value = 'Harticka'
$ws.Range("G26").Value = 'Pelloporus leporinus'
$ws.Range("H26").Value = '(Fr.) Krieglst.'
$ws.Range("Q26").Value = 447144
$ws.Range("R26").Value = 7043043
$ws.Range("S26").Value = 13
$ws.Range("Z26").Value = '15:18'
$ws.Range("AB26").Value = '15:18'

# Row 27
$ws.Range("A27").Value = 130865703
$ws.Range("B27").Value = 89194
$ws.Range("E27").Value = 510
$ws.Range("F27").Value = 'Doftskinn'
$ws.Range("G27").Value = 'Cystostereum murrayi'
$ws.Range("H27").Value = '(Berk. & M.A.Curtis.) Pouzar'
$ws.Range("Q27").Value = 447410
$ws.Range("R27").Value = 7042768
$ws.Range("S27").Value = 8
$ws.Range("Z27").Value = '14:20'
$ws.Range("AB27").Value = '14:20'

# Row 28
$ws.Range("B28").Value = 79244

# Row 30
$ws.Range("B30").Value = 91805

# Row 31
$ws.Range("B31").Value = 79244
